$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "/home/daniel/Spike Data/Matlab files/Exp 6 baseline.mat"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 10750
$ws.Range("E5").Value = 12800
$ws.Range("F5").Value = 109400
$ws.Range("G5").Value = 111200
$ws.Range("H5").Value = 189600
$ws.Range("I5").Value = 191300

# Row 6
$ws.Range("A6").Value = "/home/daniel/Spike Data/Matlab files/exp 30.mat"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 16830
$ws.Range("E6").Value = 18330
$ws.Range("F6").Value = 125400
$ws.Range("G6").Value = 127400
$ws.Range("H6").Value = 215100
$ws.Range("I6").Value = 216800

# Row 7
$ws.Range("A7").Value = "/home/daniel/Spike Data/Matlab files/exp 31.mat"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1424
$ws.Range("E7").Value = 2977
$ws.Range("F7").Value = 87050
$ws.Range("G7").Value = 89250

# Match number format used by other cells in column A (text style)
$ws.Range("A5:A7").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("H15").Select()
